$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Title paragraph: "Definición" -> "Definicion" (the run holding
#    only the accented "ó" becomes "o"), and the "_GoBack" bookmark
#    moves to sit right after that run (it used to live further down,
#    in the now-empty paragraph right under "CASOS DE PRUEBA").
# ------------------------------------------------------------------
$full = $d.Content.Text
$accentIdx = $full.IndexOf([char]0x00F3)   # the lone "ó" character

# Drop a throw-away bookmark just *before* the accented letter so this
# engine's run-merge-on-edit behaviour can't fuse it back into the
# previous "Definici" run once we rewrite its text.
$guardRange = $d.Range($accentIdx, $accentIdx)
$d.Bookmarks.Add("ZZTmpGuard", $guardRange)

# Real "_GoBack" bookmark goes right after the letter, between it and
# the following "n" run - this also keeps that run from merging in.
$goBackRange = $d.Range($accentIdx + 1, $accentIdx + 1)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# Now it's safe to flip the accented run's text without disturbing its
# neighbours.
$accentRange = $d.Range($accentIdx, $accentIdx + 1)
$accentRange.Text = "o"

# Drop the guard bookmark again now that the edit is done.
$d.Bookmarks.Item("ZZTmpGuard").Delete()

# ------------------------------------------------------------------
# 2) "CASOS DE PRUEBA" paragraph: give the paragraph mark and every
#    run in it an explicit Arial font.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "CASOS DE PRUEBA`r") {
        $pr = $para.Range
        $pr.Font.Name = "Arial"
        $pr.Font.NameBi = "Arial"
        break
    }
}
